# Users.xlsx — replace the second patient record with the new one
# (Fadi Badarni) and touch up the view state to match the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new patient record -------------------------------------------
$ws.Range("A2").Value = 209315647
$ws.Range("B2").Value = 524183083
$ws.Range("C2").Value = "fadi"
$ws.Range("D2").Value = "badarni"
$ws.Range("E2").Value = 23
$ws.Range("F2").Value = 82
$ws.Range("G2").Value = 184
$ws.Range("H2").Value = "fadybd1@gmail.com"

# Turn the e-mail address into a real mailto: hyperlink, styled the same
# way as the existing one in H1.
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:fadybd1@gmail.com")
$ws.Range("H2").Style = "Hyperlink"

# --- Column sizing ----------------------------------------------------
# The new numeric ids in columns A and B are wider, so both columns get
# resized to fit them.
$ws.Range("A1:B2").ColumnWidth = 9.1

# --- View state ---------------------------------------------------------
# Active cell/selection moved before the file was saved.
$ws.Range("O11").Select()
